$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 402; this shifts the existing rows 402..459 down to 403..460
# and extends the used range / dimension accordingly.
$ws.Rows.Item(402).Insert()

# Populate the new row 402 with the new weekly price record.
$ws.Range("A402").Value = 10
$ws.Range("B402").Value = "Vega Modelo de Temuco"
$ws.Range("C402").Value = "La Araucanía"
$ws.Range("D402").Value = 45124
$ws.Range("E402").Value = 9
$ws.Range("F402").Value = "Fruta"
$ws.Range("G402").Value = 100102
$ws.Range("H402").Value = "Cítricos"
$ws.Range("I402").Value = 100102006
$ws.Range("J402").Value = "Pomelo"
$ws.Range("K402").Value = "Start Ruby"
$ws.Range("L402").Value = "Primera"
$ws.Range("M402").Value = 55
$ws.Range("N402").Value = 15000
$ws.Range("O402").Value = 15000
$ws.Range("P402").Value = 15000
$ws.Range("Q402").Value = '$/bandeja 15 kilos granel'
$ws.Range("R402").Value = "Región de O'Higgins"
$ws.Range("S402").Value = 1000
$ws.Range("T402").Value = 15
